$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 6) continuing the "Recursion" topic column,
# mirroring the pattern already present in rows 2-5.
$ws.Range("A6").Value = "Recursion"

# After entering the value, the active cell in Excel moves to the right,
# leaving B6 selected (matching the saved selection in the sheet).
$ws.Range("B6").Select()
